$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Changes" text for the 0.1.0 release (row 3, column C)
$ws.Range("C3").Value = "Nearly all evaluations from Cakedefi-Review.com integrated`nAdaptions made for responsive layout, especially sidebar"

# Row 3 height shrinks from 45 to 30 now that the text is shorter
$ws.Rows(3).RowHeight = 30

# Update the active selection to C4 only
$null = $ws.Range("C4").Select()
